$d = $word.ActiveDocument

# The document ends with an empty, bulleted (ListParagraph/numPr) paragraph
# that currently has no runs at all, followed by a trailing (non-bulleted)
# blank paragraph. We need to:
#   1. Give that empty bulleted paragraph the text "User Operations Screen".
#   2. Insert a brand-new bulleted paragraph right after it containing the
#      text "User Menu Screen" followed by a manual line break.

# Locate the empty bulleted paragraph (the last ListParagraph/numPr
# paragraph that has no text yet) by scanning from the end.
$targetIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -eq 0 -and $p.Range.ListFormat.ListType -ne 0) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq -1) {
    $targetIndex = $d.Paragraphs.Count - 1
}

$target = $d.Paragraphs.Item($targetIndex)

# 1) Add "User Operations Screen" text into the existing empty paragraph.
$target.Range.InsertBefore("User Operations Screen")
$d.Paragraphs.Item($targetIndex).Range.LanguageID = "en-US"

# 2) Insert a brand new paragraph right after it, then overwrite that new
#    paragraph's whole range (run content) with the exact OOXML we need so
#    the trailing manual line-break run keeps its own <w:rPr>.
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' + `
    '<w:rPr><w:lang w:val="en-US"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>User Menu Screen</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:br/></w:r>' + `
'</w:p>'

$newParaRange.InsertXML($newParaXml)

Write-Output "done"
